$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as literal text in the source data
# (e.g. "0.550", "10.50") so we force Text format before assigning, then
# restore the default "Normal" style afterwards so no stray number format
# is left applied to the cell (matches original: no explicit style on data cells).
$dCells = @('D2','D3','D5','D6','D7','D8','D9','D10','D12','D13','D14','D15','D17','D18','D19','D20','D21','D22','D23','D24','D25','D26','D27','D30','D32','D33','D34','D36','D37','D39','D40','D43','D44','D46','D47','D48','D49','D50','D51')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.406.57"
$ws.Range("E2").Value = "  +5.19%  "
$ws.Range("D3").Value = "3.181.77"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("D5").Value = "401.55"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").Value = "108.77"
$ws.Range("E6").Value = "  +4.93%  "
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("D10").Value = "38.89"
$ws.Range("E10").Value = "  +4.86%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "0.0881"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").Value = "3.671.15"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "19.06"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "8.03"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("E16").Value = "  +8.89%  "
$ws.Range("D17").Value = "3.180.31"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "10.50"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "54.359.87"
$ws.Range("E19").Value = "  +4.88%  "
$ws.Range("D20").Value = "3.32"
$ws.Range("E20").Value = "  +3.90%  "
$ws.Range("D21").Value = "12.87"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").Value = "0.0₃0991"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "71.94"
$ws.Range("E23").Value = "  +2.70%  "
$ws.Range("D24").Value = "274.32"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").Value = "3.27"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").Value = "8.02"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "27.71"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +4.07%  "
$ws.Range("D32").Value = "11.07"
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("D33").Value = "0.0497"
$ws.Range("E33").Value = "  +10.20%  "
$ws.Range("D34").Value = "37.01"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "50.86"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  +7.47%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  +9.12%  "
$ws.Range("D40").Value = "4.16"
$ws.Range("E40").Value = "  +12.44%  "
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").Value = "17.28"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").Value = "130.02"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").Value = "22.32"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").Value = "2.06"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "2.091.00"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").Value = "0.0345"
$ws.Range("E50").Value = "  +8.28%  "
$ws.Range("D51").Value = "0.0508"
$ws.Range("E51").Value = "  +11.38%  "

# Restore default styling on the D cells we touched (removes the temporary
# Text number-format so the cell goes back to having no explicit style,
# exactly like the untouched data cells).
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
